# svx/qa/unit/data/video-snapshot.pptx: add crop + move/resize to the
# media snapshot picture, and bump the cached "today" date placeholder
# text (8/23/2022 -> 8/25/2022) on the slide master and every slide
# layout.

$p = $ppt.ActivePresentation

# --- 1) Crop the media snapshot picture and reposition/resize it -----
# Original: off (3048000,1143000) ext (6096000,4572000), no crop.
# Target:   off (4661452,1143000) ext (2991678,4572000),
#           srcRect l="25000" r="25000" (25% cropped off each side).
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)

$shp.PictureFormat.CropLeft = 120
$shp.PictureFormat.CropRight = 120

$shp.Left = 367.0435
$shp.Width = 235.5651968503937

# --- 2) Update the cached date placeholder text everywhere it appears -
function Update-DatePlaceholder($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $cand = $shapes.Item($i)
        if ($cand.HasTextFrame -and $cand.PlaceholderFormat.Type -eq 16) {
            if ($cand.TextFrame.TextRange.Text -eq "8/23/2022") {
                $cand.TextFrame.TextRange.Text = "8/25/2022"
            }
        }
    }
}

$master = $s.Master
Update-DatePlaceholder($master)

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder($layout)
}
